# StructureDefinition-claim-response-benefit-plan.xlsx
# Rebrand IBM/Alvearie -> LinuxForHealth, bump version/date, and fix a
# duplicated FHIR constraint that incorrectly appeared on the top-level
# "Extension" row instead of only on "Extension.extension".

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL (row 2 / B2): ibm.com -> linuxforhealth.org
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-response-benefit-plan"

# Version (row 3 / B3): 7.0.0 -> 8.0.0
$wsMeta.Range("B3").Value = "8.0.0"

# Date (row 8 / B8)
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (row 9 / B9): Alvearie Team -> LinuxForHealth Team
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension") incorrectly carried the ele-1/ext-1 constraint that
# belongs only to row 4 ("Extension.extension"); clear it here.
$wsElem.Range("AI2").Value = ""

# The "Fixed Value" on Extension.url mirrors the StructureDefinition URL.
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-response-benefit-plan"
